# poisson_naive versao media ponderada
# Update the "index"/weight column (A) for each match row so that it
# reflects the new weighted-average based numbering instead of the
# previous simple sequential numbering.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 313
    3  = 315
    4  = 316
    5  = 318
    6  = 320
    7  = 322
    8  = 323
    9  = 326
    10 = 327
    11 = 329
    12 = 331
    13 = 333
    14 = 336
    15 = 28
    16 = 71
    17 = 96
    18 = 111
    19 = 179
    20 = 198
    21 = 285
    22 = 298
    23 = 338
    24 = 386
    25 = 399
    26 = 429
    27 = 499
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}
